# Weekly data refresh: insert this week's new price row at the top of the
# data block (row 43) and push the previously-existing rows (old 43..68)
# down by one (new 44..69), matching the source feed's newest-first order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 43; Excel shifts row 43..68 -> 44..69
# and copies formatting (e.g. the date-format style on column D) down from
# the row above, same as a manual "Insert Sheet Rows" in the UI.
$ws.Rows.Item(43).Insert()

# Populate the freshly inserted row 43 with the new week's observation.
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C43").Value = "Arica y Parinacota"
$ws.Range("D43").Value = 45086
$ws.Range("E43").Value = 15
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100108
$ws.Range("H43").Value = "Tropicales y subtropicales"
$ws.Range("I43").Value = 100108001
$ws.Range("J43").Value = "Guayaba"
$ws.Range("K43").Value = "Sin especificar"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 430
$ws.Range("N43").Value = 4000
$ws.Range("O43").Value = 7000
$ws.Range("P43").Value = 5698
$ws.Range("Q43").Value = "$/caja 10 kilos"
$ws.Range("R43").Value = "Región de Arica y Parinacota"
$ws.Range("S43").Value = 570
$ws.Range("T43").Value = 10
